$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D:D").Insert()

# Copy number styles from column E (which now holds old column D's formatting) into D
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)  # xlPasteFormats = -4122
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

